$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 15 (the "period" extension slice): fill in the
#     specific extension Type(s)/Short/Definition instead of the
#     generic "Extension" placeholder values.
$ws.Range("J15").Value = "Extension {http://hl7.org/fhir/ca-bc/provider/StructureDefinition/bc-period-extension}`n"
$ws.Range("K15").Value = "BC Valid Period"
$ws.Range("L15").Value = "The period for when the extended element is valid."

# --- Remove the four sub-element rows that used to spell out
#     Extension.extension.id / .extension / .url / .value[x] for the
#     "period" slice (rows 16-19). This pulls the old rows 20-22 up
#     to become the new rows 16-18.
$ws.Rows("16:19").Delete()

# --- Rename the "status" extension slice (now at row 16) to
#     "endReason" and point it at the new extension definition.
$ws.Range("B16").Value = "endReason"
$ws.Range("J16").Value = "Extension {http://hl7.org/fhir/ca-bc/provider/StructureDefinition/bc-end-reason-extension}`n"
$ws.Range("K16").Value = "BC End Reason Extension"
$ws.Range("L16").Value = "Tracking end reasons."

# --- Restore row visibility. Deleting rows above cleared every row's
#     Hidden flag in this runtime, so re-apply the full hidden set for
#     the (post-delete) sheet.
$ws.Rows(1).Hidden = $false
$ws.Rows(2).Hidden = $true
$ws.Rows(3).Hidden = $true
$ws.Rows(4).Hidden = $true
$ws.Rows(5).Hidden = $false
$ws.Rows(6).Hidden = $true
$ws.Rows(7).Hidden = $true
$ws.Rows(8).Hidden = $true
$ws.Rows(9).Hidden = $true
$ws.Rows(10).Hidden = $false
$ws.Rows(11).Hidden = $true
$ws.Rows(12).Hidden = $true
$ws.Rows(13).Hidden = $true
$ws.Rows(14).Hidden = $true
$ws.Rows(15).Hidden = $false
$ws.Rows(16).Hidden = $false
$ws.Rows(17).Hidden = $true
$ws.Rows(18).Hidden = $true

# --- Re-point the AutoFilter at the new, smaller range (A1:AJ18) and
#     re-apply the two custom filters that were on it before.
$ws.AutoFilterMode = $false
$ws.Range("A1:AJ18").AutoFilter(7, "<>" + " ")
$ws.Range("A1:AJ18").AutoFilter(27, @(""), 7)

# --- Conditional formatting previously covered A2:AI21; shrink it to
#     A2:AI17 to track the smaller sheet (keeps both existing rules &
#     their dxf links intact).
$fcs = $ws.Range("A2:AI21").FormatConditions
$fcs.Item(1).ModifyAppliesToRange($ws.Range("A2:AI17"))

# --- The hidden _xlnm._FilterDatabase defined name also needs to
#     track the new, smaller used range.
$wb.Names.Item(1).RefersTo = "=Elements!`$A`$1:`$AJ`$18"
